$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q: header + weight/card data -------------------------------

# Header cell Q3, formatted like the other header cells (P3)
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q3").Value = "权重/单卡alllayers"

# Data cells Q4:Q11, formatted like the neighboring numeric column C4:C11
$ws.Range("C4:C11").Copy()
$ws.Range("Q4:Q11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("Q4").Value = 838860800
$ws.Range("Q5").Value = 671088640
$ws.Range("Q6").Value = 4194304000
$ws.Range("Q7").Value = 2097152000
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0

# Match column width of the other "12"-wide columns (e.g. K)
$ws.Columns.Item(17).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# --- New summary row 25 -----------------------------------------------------

$ws.Range("A15").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A25").Value = "权重显存/单卡(GB)"

$ws.Range("B15").Copy()
$ws.Range("B25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B25").Value = 7.266

# --- Merged title cell now spans through the new column --------------------

$ws.Range("A1:O1").UnMerge()
$ws.Range("A1:Q1").Merge()
